# Add a new row to the end of the (single) summary table in the report,
# describing the newly implemented "Calculate fee amount" feature.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Append a new row after the last existing row ("Get list of overdue
# tickets ..."); Word clones the formatting (fonts, centered % column)
# from the row it is appended after.
$newRow = $t.Rows.Add()

$cells = $newRow.Cells
$cells.Item(1).Range.Text = "Calculate fee amount "
$cells.Item(2).Range.Text = "100"
$cells.Item(3).Range.Text = "Lost books penalty + overdue fee"
